# Daily attendance processing - 2025-12-19 10:58:45
# Normalize the "Recorded By" column (G) so that every comma-separated
# list of recorders is rotated to start at the literal "System" entry
# (the rest of the entries keep their original relative order, wrapping
# around). Rows whose value does not contain an exact "System" token are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $parts = $val -split ", "
    $idx = [System.Array]::IndexOf($parts, "System")

    if ($idx -gt 0) {
        $rotated = $parts[$idx..($parts.Count - 1)] + $parts[0..($idx - 1)]
        $newVal = $rotated -join ", "
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
